$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6) need their Fecha/Volumen/Precio.../Origen values
# cyclically shifted up by one row: row N takes the values that were
# previously on row N+1, and row 6 takes the values that were on row 2.

$rows = @(2, 3, 4, 5, 6)

# Capture original values (D, M, N, O, P, R, S) for each row before any writes.
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        D = $ws.Range("D$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        R = $ws.Range("R$r").Value2
        S = $ws.Range("S$r").Value2
    }
}

# Map: new value for row r = original value of the "source" row.
$sourceRow = @{ 2 = 3; 3 = 4; 4 = 5; 5 = 6; 6 = 2 }

foreach ($r in $rows) {
    $src = $sourceRow[$r]
    $vals = $orig[$src]

    $ws.Range("D$r").Value2 = $vals.D
    $ws.Range("M$r").Value2 = $vals.M
    $ws.Range("N$r").Value2 = $vals.N
    $ws.Range("O$r").Value2 = $vals.O
    $ws.Range("P$r").Value2 = $vals.P
    $ws.Range("R$r").Value2 = $vals.R
    $ws.Range("S$r").Value2 = $vals.S
}
